# This script swaps the values of columns B, D, E, F, G between each of the
# following pairs of adjacent rows. Columns A (serial no.) and C (item name)
# stay with their original row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowPairs = @(
    @(136, 137),
    @(163, 164),
    @(292, 293),
    @(294, 295),
    @(311, 312),
    @(315, 316),
    @(356, 357),
    @(420, 421),
    @(465, 466),
    @(472, 473),
    @(490, 491),
    @(596, 597),
    @(705, 706),
    @(732, 733)
)

$cols = @("B", "D", "E", "F", "G")

foreach ($pair in $rowPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    foreach ($col in $cols) {
        $addr1 = "$col$r1"
        $addr2 = "$col$r2"

        $v1 = $ws.Range($addr1).Value2
        $v2 = $ws.Range($addr2).Value2

        $ws.Range($addr1).Value2 = $v2
        $ws.Range($addr2).Value2 = $v1
    }
}
